# Refs #435. Fixed an issue to write SampleToSi log.
#
# Slide "Test Case" (the slide that documents Test 1 / Test 3 sample &
# detector positions) needs two text edits in the Content Placeholder:
#
#   1. "Sample is located at (0, -0.088, 0)"
#        -> "Sample is located at (0, 0, 0-0.088)"
#      (the y/z values were transposed - this is split into two runs in
#      the authored OOXML, so we recreate that by rewriting the tail of
#      the paragraph text, which causes a run split.)
#
#   2. "Sample detector distance is then 10 m"
#        -> the literal value "10" is called out in red, matching the
#      other "Sample detector distance is ..." paragraphs on the same
#      slide that already highlight their numeric value in red.

$p = $ppt.ActivePresentation

$targetSlide = $null
$targetShape = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            $fullText = $shape.TextFrame.TextRange.Text
            if ($fullText.Contains("Sample is located at (0, -0.088, 0)") -and $fullText.Contains("Sample detector distance is then 10 m")) {
                $targetSlide = $slide
                $targetShape = $shape
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange
$paraCount = $tr.Paragraphs().Count

# --- Edit 1: "Sample is located at (0, -0.088, 0)" ---------------------
for ($i = 1; $i -le $paraCount; $i++) {
    $para = $tr.Paragraphs($i, 1)
    $t = $para.Text
    if ($t.StartsWith("Sample is located at (0, -0.088, 0)")) {
        $prefix = "Sample is located at (0"
        $suffixStart = $prefix.Length + 1
        $suffixLen = "Sample is located at (0, -0.088, 0)".Length - $prefix.Length
        $tail = $para.Characters($suffixStart, $suffixLen)
        $tail.Text = ", 0, 0-0.088)"
    }
}

# --- Edit 2: "Sample detector distance is then 10 m" -------------------
for ($i = 1; $i -le $paraCount; $i++) {
    $para = $tr.Paragraphs($i, 1)
    $t = $para.Text
    if ($t.StartsWith("Sample detector distance is then 10 m")) {
        $numStart = "Sample detector distance is then ".Length + 1
        $numLen = "10".Length
        $numRange = $para.Characters($numStart, $numLen)
        $numRange.Font.Color.RGB = 255
    }
}
